$wb = $excel.ActiveWorkbook

# --- Overview sheet: update status for the second file (e6052173-...md) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: update status + latest handoff datetime for row 3 ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "2016-02-24 07:09:40"

# --- de-de sheet: update status + latest handoff datetime for row 3 ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "2016-02-24 07:09:52"
